# All data added to database.
# Remove the "Undergraduate transfer-in enrollment" column (E) and the
# "First-time, full-time bachelor's seeking student retention rate" column
# (originally H, becomes G after the first delete) from the Sheet1 stats
# table. This shifts the old Student-to-faculty ratio column into E and the
# old Percent-of-all-students-who-are-female column into F.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(5).Delete() | Out-Null
$ws.Columns.Item(7).Delete() | Out-Null

# Relabel the surviving header cells to match the new column meanings.
$ws.Cells.Item(1, 6).Value2 = "Percent of female student"
$ws.Cells.Item(1, 5).Value2 = "Student_faculty ratio"

# The header row's autofit height shrinks now that it only needs to wrap
# two (shorter) headers instead of four.
$ws.Rows.Item(1).RowHeight = 17

# Reset the sheet's selection back to the top of the (now narrower) table.
$ws.Range("E1").Select() | Out-Null
